# Avances Etiquetado Roboflow - Miercoles 5/28/2025
# Fill in the missed Tuesday (27/5/2025) and Wednesday (28/5/2025) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 used to be the "next empty templated row" right below the table's
# last data row (old row 24). Copy its current formatting down to row 26
# before row 24 gets real data, so row 26 ends up with the same blank
# template look the table always keeps right under the last filled row.
$ws.Range("D24:J24").Copy()
$ws.Range("D26:J26").PasteSpecial(-4122)
$ws.Range("D26:J26").ClearContents()

# Row 25 becomes a real data row too - give it the same formatting as the
# preceding data rows (e.g. row 23) before filling in values.
$ws.Range("D23:J23").Copy()
$ws.Range("D25:J25").PasteSpecial(-4122)

# Tuesday 27/5/2025 (row 24)
$ws.Range("D24").Value = "27/5/2025"
$ws.Range("E24").Value = 125
$ws.Range("F24").Value = 441
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 650
$ws.Range("J24").Value = "N/A"

# Wednesday 28/5/2025 (row 25)
$ws.Range("D25").Value = "28/5/2025"
$ws.Range("E25").Value = 110
$ws.Range("F25").Value = 456
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 650
$ws.Range("J25").Value = "N/A"

# Table1 grows by two rows to keep including the new data.
$wb.Worksheets.Item("Sheet1").ListObjects("Table1").Resize($ws.Range("D4:J26")) | Out-Null

# Update the saved selection the same way Excel would after this edit.
$ws.Range("H29").Select() | Out-Null
